# False cobT/cobS genes are members of B12 pathway
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_pathway_genes")

# Rows 68 & 69 (false_S -> cobS-related): chel_gene "false_cobS" first
$ws.Range("E68").Value = "false_cobS"
$ws.Range("E69").Value = "false_cobS"

# Rows 60 & 61 (false_M -> cobT-related): chel_gene_group "cobT_cobU"
$ws.Range("D60").Value = "cobT_cobU"
$ws.Range("D61").Value = "cobT_cobU"

# Rows 60 & 61: chel_gene "false_cobT"
$ws.Range("E60").Value = "false_cobT"
$ws.Range("E61").Value = "false_cobT"

# Rows 68 & 69: chel_gene_group "cobV_cobS" (reuses existing shared string)
$ws.Range("D68").Value = "cobV_cobS"
$ws.Range("D69").Value = "cobV_cobS"

# All four rows: chel_pathway becomes "B12", with explicit black font color
$ws.Range("C60").Value = "B12"
$ws.Range("C60").Font.Color = 0
$ws.Range("C61").Value = "B12"
$ws.Range("C61").Font.Color = 0
$ws.Range("C68").Value = "B12"
$ws.Range("C68").Font.Color = 0
$ws.Range("C69").Value = "B12"
$ws.Range("C69").Font.Color = 0

# View change: move the active selection to D9
$ws.Activate()
$ws.Range("D9").Select()
